$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# The "Periodo Mora" table (rows 16-27, column E = period, column F = Valor Mora) was sorted
# in descending order (1911 down to 1812). The database update re-sorts these same twelve
# periods in ascending order (1812 up to 1911), carrying each period's "Valor Mora" amount
# along with it, and folds in the new part-1 data for the account statement.

$periods = @("1812", "1901", "1902", "1903", "1904", "1905", "1906", "1907", "1908", "1909", "1910", "1911")
$valores = @(26041, 31249, 31249, 31249, 31249, 31249, 31249, 31249, 31249, 31249, 31249, 28124)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $valores[$i]
}
